{"js": "// \"Added name to Project Plan\"\n// The placeholder paragraph \"Student Names\" (on the title page, directly\n// below the \"<Project Name>\" line) is replaced with the actual student's\n// name, \"Brianne Byer\".\n\nconst body = context.document.body;\n\n// Locate the placeholder text. Using search() (rather than indexing into\n// body.paragraphs) is robust even though \"Student Names\" is split across\n// two separate runs (\"Student Name\" + \"s\") in the original document.\nconst results = body.search(\"Student Names\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole matched range (both runs) with a single new run\n  // containing the student's name.\n  results.items[0].insertText(\"Brianne Byer\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# \"Added name to Project Plan\"\n# The placeholder paragraph \"Student Names\" (on the title page, directly\n# below the \"<Project Name>\" line) is replaced with the actual student's\n# name, \"Brianne Byer\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Student Names\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Brianne Byer\"\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n# Replace(wdReplaceAll=2) - swap every occurrence (there is exactly one) in\n# a single pass.\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
